$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the frequency-offset values in row 2 (previously placeholder 0s) ---
$ws.Range("A2").Value = 4961272.0199999996
$ws.Range("C2").Value = 14866726.869999999
$ws.Range("E2").Value = 24774372.23
$ws.Range("G2").Value = 34680641.270000003
$ws.Range("I2").Value = 44587621.439999998
$ws.Range("K2").Value = 54494404.460000001
$ws.Range("M2").Value = 64401754.57

# --- Header row (A1:N1) goes back to the workbook's default (unstyled) look ---
$ws.Range("A1:N1").Style = "Normal"

# --- New row heights for the header + data row ---
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 15.75

# --- Column E is now sized to fit its (much wider) new content ---
$ws.Columns.Item(5).ColumnWidth = 10.75

# --- Give the updated frequency cells a medium box border, centered/wrapped text ---
$freqRng = $ws.Range("A2,C2,E2,G2,I2,K2,M2")
$freqRng.Borders.Weight = -4138
$freqRng.VerticalAlignment = -4108
$freqRng.WrapText = $true

# --- Selection ends up on M2 ---
$ws.Range("M2").Select()
